$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)
$ws.Range("B7").Value = "test"
